$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.972.52'
$ws.Cells.Item(2, 5).Value = '  +0.52%  '
$ws.Cells.Item(3, 4).Value = '1.894.38'
$ws.Cells.Item(3, 5).Value = '  +0.15%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9998'
$ws.Cells.Item(4, 4).NumberFormat = 'General'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '0.7753'
$ws.Cells.Item(5, 4).NumberFormat = 'General'
$ws.Cells.Item(5, 5).Value = '  -0.36%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '243.96'
$ws.Cells.Item(6, 4).NumberFormat = 'General'
$ws.Cells.Item(6, 5).Value = '  +0.18%  '
$ws.Cells.Item(7, 5).Value = '  -0.15%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3132'
$ws.Cells.Item(8, 4).NumberFormat = 'General'
$ws.Cells.Item(8, 5).Value = '  +0.21%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '25.81'
$ws.Cells.Item(9, 4).NumberFormat = 'General'
$ws.Cells.Item(9, 5).Value = '  +2.27%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.07259'
$ws.Cells.Item(10, 4).NumberFormat = 'General'
$ws.Cells.Item(10, 5).Value = '  +1.65%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.08681'
$ws.Cells.Item(11, 4).NumberFormat = 'General'
$ws.Cells.Item(11, 5).Value = '  +7.47%  '
$ws.Cells.Item(12, 4).Value = '2.031.32'
$ws.Cells.Item(12, 5).Value = '  +5.98%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.7730'
$ws.Cells.Item(13, 4).NumberFormat = 'General'
$ws.Cells.Item(13, 5).Value = '  +1.67%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '5.422'
$ws.Cells.Item(14, 4).NumberFormat = 'General'
$ws.Cells.Item(14, 5).Value = '  -0.48%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '94.55'
$ws.Cells.Item(15, 4).NumberFormat = 'General'
$ws.Cells.Item(15, 5).Value = '  +2.69%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '6.226'
$ws.Cells.Item(16, 4).NumberFormat = 'General'
$ws.Cells.Item(16, 5).Value = '  +1.66%  '
$ws.Cells.Item(17, 4).Value = '30.072.28'
$ws.Cells.Item(17, 5).Value = '  +1.01%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '13.95'
$ws.Cells.Item(18, 4).NumberFormat = 'General'
$ws.Cells.Item(18, 5).Value = '  +0.28%  '
$ws.Cells.Item(19, 4).Value = '2.308.70'
$ws.Cells.Item(19, 5).Value = '  +13.99%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '245.94'
$ws.Cells.Item(20, 4).NumberFormat = 'General'
$ws.Cells.Item(20, 5).Value = '  +1.22%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.000007895'
$ws.Cells.Item(21, 4).NumberFormat = 'General'
$ws.Cells.Item(21, 5).Value = '  +1.76%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '8.201'
$ws.Cells.Item(22, 4).NumberFormat = 'General'
$ws.Cells.Item(22, 5).Value = '  +1.69%  '
$ws.Cells.Item(23, 5).Value = '  -0.05%  '
$ws.Cells.Item(24, 5).Value = '  -0.10%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.1601'
$ws.Cells.Item(25, 4).NumberFormat = 'General'
$ws.Cells.Item(25, 5).Value = '  -1.23%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '9.542'
$ws.Cells.Item(26, 4).NumberFormat = 'General'
$ws.Cells.Item(26, 5).Value = '  +1.70%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '163.08'
$ws.Cells.Item(27, 4).NumberFormat = 'General'
$ws.Cells.Item(27, 5).Value = '  +0.59%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '18.87'
$ws.Cells.Item(28, 4).NumberFormat = 'General'
$ws.Cells.Item(28, 5).Value = '  +1.05%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.050'
$ws.Cells.Item(29, 4).NumberFormat = 'General'
$ws.Cells.Item(29, 5).Value = '  +0.43%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.432'
$ws.Cells.Item(30, 4).NumberFormat = 'General'
$ws.Cells.Item(30, 5).Value = '  +1.47%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.546'
$ws.Cells.Item(31, 4).NumberFormat = 'General'
$ws.Cells.Item(31, 5).Value = '  +0.04%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '4.526'
$ws.Cells.Item(32, 4).NumberFormat = 'General'
$ws.Cells.Item(32, 5).Value = '  +1.30%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '4.132'
$ws.Cells.Item(33, 4).NumberFormat = 'General'
$ws.Cells.Item(33, 5).Value = '  +0.82%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.05453'
$ws.Cells.Item(34, 4).NumberFormat = 'General'
$ws.Cells.Item(34, 5).Value = '  -1.20%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.252'
$ws.Cells.Item(35, 4).NumberFormat = 'General'
$ws.Cells.Item(35, 5).Value = '  -0.75%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.7547'
$ws.Cells.Item(36, 4).NumberFormat = 'General'
$ws.Cells.Item(36, 5).Value = '  +1.85%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.002'
$ws.Cells.Item(37, 4).NumberFormat = 'General'
$ws.Cells.Item(37, 5).Value = '  +0.91%  '
$ws.Cells.Item(38, 5).Value = '  +2.62%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01979'
$ws.Cells.Item(39, 4).NumberFormat = 'General'
$ws.Cells.Item(39, 5).Value = '  +3.62%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.786'
$ws.Cells.Item(40, 4).NumberFormat = 'General'
$ws.Cells.Item(40, 5).Value = '  +0.11%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.4527'
$ws.Cells.Item(41, 4).NumberFormat = 'General'
$ws.Cells.Item(41, 5).Value = '  +2.86%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '74.01'
$ws.Cells.Item(42, 4).NumberFormat = 'General'
$ws.Cells.Item(42, 5).Value = '  +0.57%  '
$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Value = '1.094.79'
$ws.Cells.Item(43, 5).Value = '  -4.00%  '
$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '6.066'
$ws.Cells.Item(44, 4).NumberFormat = 'General'
$ws.Cells.Item(44, 5).Value = '  +3.84%  '
$ws.Cells.Item(45, 4).Value = '2.229.07'
$ws.Cells.Item(45, 5).Value = '  +12.02%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.8540'
$ws.Cells.Item(46, 4).NumberFormat = 'General'
$ws.Cells.Item(46, 5).Value = '  +0.46%  '
$ws.Cells.Item(47, 5).Value = '  -0.16%  '
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '103.14'
$ws.Cells.Item(48, 4).NumberFormat = 'General'
$ws.Cells.Item(48, 5).Value = '  -0.20%  '
$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.888'
$ws.Cells.Item(49, 4).NumberFormat = 'General'
$ws.Cells.Item(49, 5).Value = '  +1.20%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '7.636'
$ws.Cells.Item(50, 4).NumberFormat = 'General'
$ws.Cells.Item(50, 5).Value = '  +2.82%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '9.925'
$ws.Cells.Item(51, 4).NumberFormat = 'General'
$ws.Cells.Item(51, 5).Value = '  +0.12%  '
